$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("J19").Value = "fini"
$ws.Range("C21").Value = "other_a"
$ws.Range("D21").Value = "other_b"
$ws.Range("E21").Value = "other_c"
$ws.Range("C22").Value = 8
$ws.Range("D22").Value = "horse strike equal"
$ws.Range("E22").Value = "soil cells inch"
$ws.Range("C23").Value = 13
$ws.Range("D23").Value = "health region away"
$ws.Range("E23").Value = "scene still wire"
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = "remember change fair"
$ws.Range("E24").Value = "prepare spell left"
[void]$ws.Range("J20").Select()
